# Apply "new card stack (balanced) - 104 cards total" edit:
# Update the "copies" column (column G) values on the three market sheets.

$wb = $excel.ActiveWorkbook

# --- RocketMarket sheet ---
$ws = $wb.Worksheets.Item("RocketMarket")
$ws.Range("G2").Value = 2
$ws.Range("G3").Value = 2
$ws.Range("G4").Value = 2
$ws.Range("G5").Value = 2
$ws.Range("G8").Value = 2
$ws.Range("G9").Value = 1
$ws.Range("G10").Value = 1

# --- ShieldMarket sheet ---
$ws = $wb.Worksheets.Item("ShieldMarket")
$ws.Range("G2").Value = 2
$ws.Range("G3").Value = 2
$ws.Range("G4").Value = 2
$ws.Range("G5").Value = 2
$ws.Range("G7").Value = 1
$ws.Range("G8").Value = 1

# --- SpecialsMarket sheet ---
$ws = $wb.Worksheets.Item("SpecialsMarket")
$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 1
$ws.Range("G6").Value = 1
$ws.Range("G7").Value = 1
$ws.Range("G8").Value = 2
$ws.Range("G9").Value = 2
